$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting existing rows 304-418 down to 305-419.
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new Brocoli record.
$ws.Range("A304").Value = 5
$ws.Range("B304").Value = "Macroferia Regional de Talca"
$ws.Range("C304").Value = "Maule"
$ws.Range("D304").Value = 44825
$ws.Range("E304").Value = 7
$ws.Range("F304").Value = 100112023
$ws.Range("G304").Value = "Brócoli"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 3000
$ws.Range("K304").Value = 1000
$ws.Range("L304").Value = 1000
$ws.Range("M304").Value = 1000
$ws.Range("N304").Value = "$/unidad"
$ws.Range("O304").Value = "Región del Maule"
$ws.Range("P304").Value = 1000
$ws.Range("Q304").Value = 1
$ws.Range("R304").Value = "Hortaliza"
